$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks (L and N widen) ---
$ws.Columns(12).ColumnWidth = 20.166666666666668
$ws.Columns(14).ColumnWidth = 15.42578125

# --- New "TipoUsuario" / "Usuario" tables, rows 28-31 ---

# Row 28: section headers
$ws.Range("H28").Value = "TipoUsuario"
$ws.Range("H28").Interior.Color = 192
$ws.Range("H28").HorizontalAlignment = -4108
$ws.Range("I28").Interior.Color = 192
$ws.Range("I28").HorizontalAlignment = -4108
$ws.Range("H28:I28").Merge()

$ws.Range("K28").Interior.Color = 10498160

$ws.Range("L28").Value = "Usuario"
$ws.Range("L28").Interior.Color = 10498160

$ws.Range("M28").Interior.Color = 10498160
$ws.Range("M28").HorizontalAlignment = -4108
$ws.Range("N28").Interior.Color = 10498160
$ws.Range("N28").HorizontalAlignment = -4108
$ws.Range("M28:N28").Merge()

# Row 29: column headers for both tables
$ws.Range("H29").Value = "idTipoUsuario"
$ws.Range("H29").Interior.Color = 8750575
$ws.Range("H29").HorizontalAlignment = -4108

$ws.Range("I29").Value = "titulo"
$ws.Range("I29").Interior.Color = 8750575
$ws.Range("I29").HorizontalAlignment = -4108

$ws.Range("K29").Value = "idUsuario"
$ws.Range("K29").Interior.Color = 15718895
$ws.Range("K29").HorizontalAlignment = -4108

$ws.Range("L29").Value = "email"
$ws.Range("L29").Interior.Color = 15718895
$ws.Range("L29").HorizontalAlignment = -4108

$ws.Range("M29").Value = "senha"
$ws.Range("M29").Interior.Color = 15718895
$ws.Range("M29").HorizontalAlignment = -4108

$ws.Range("N29").Value = "idTipoUsuario"
$ws.Range("N29").Interior.Color = 15718895
$ws.Range("N29").HorizontalAlignment = -4108

# Row 30: Administrador / admin
$ws.Range("H30").Value = 1
$ws.Range("H30").Interior.Color = 8750575
$ws.Range("H30").HorizontalAlignment = -4108

$ws.Range("I30").Value = "Administrador"
$ws.Range("I30").Interior.Color = 8750575
$ws.Range("I30").HorizontalAlignment = -4108

$ws.Range("K30").Value = 1
$ws.Range("K30").Interior.Color = 15718895
$ws.Range("K30").HorizontalAlignment = -4108

$ws.Range("L30").Value = "admin@admin.com"
$ws.Range("L30").Interior.Color = 15718895
$ws.Range("L30").HorizontalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("L30"), "mailto:admin@admin.com")

$ws.Range("M30").Value = "admin"
$ws.Range("M30").Interior.Color = 15718895
$ws.Range("M30").HorizontalAlignment = -4108

$ws.Range("N30").Value = 1
$ws.Range("N30").Interior.Color = 15718895
$ws.Range("N30").HorizontalAlignment = -4108

# Row 31: Jogador / jogador
$ws.Range("H31").Value = 2
$ws.Range("H31").Interior.Color = 8750575
$ws.Range("H31").HorizontalAlignment = -4108

$ws.Range("I31").Value = "Jogador"
$ws.Range("I31").Interior.Color = 8750575
$ws.Range("I31").HorizontalAlignment = -4108

$ws.Range("K31").Value = 2
$ws.Range("K31").Interior.Color = 15718895
$ws.Range("K31").HorizontalAlignment = -4108

$ws.Range("L31").Value = "jogador@jogador.com"
$ws.Range("L31").Interior.Color = 15718895
$ws.Range("L31").HorizontalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("L31"), "mailto:jogador@jogador.com")

$ws.Range("M31").Value = "jogador"
$ws.Range("M31").Interior.Color = 15718895
$ws.Range("M31").HorizontalAlignment = -4108

$ws.Range("N31").Value = 2
$ws.Range("N31").Interior.Color = 15718895
$ws.Range("N31").HorizontalAlignment = -4108

# --- View settings ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("L28").Select()
